$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column G holds the "calibration_file" values; rename the codes used
# in rows 2-7 to their more descriptive equivalents.
$ws.Range("G2").Value = "calibration"
$ws.Range("G3").Value = "calibration"
$ws.Range("G4").Value = "deriv_calibration"
$ws.Range("G5").Value = "deriv_calibration"
$ws.Range("G6").Value = "calibration"
$ws.Range("G7").Value = "calibration"
